$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.809159398078918
$ws.Range("B1").Value = 2.590580224990845
$ws.Range("C1").Value = 2.70127010345459
$ws.Range("D1").Value = 3.080679178237915
$ws.Range("E1").Value = 3.341838121414185
